# CSU05 - Inscrever na corrida
# "Casos De Uso Corrigido e Nome de Telas"
#
# Renames the screen/prototype references used throughout the use-case
# document and fixes the wording of the "back to home" action text.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$Replace
    )
    $d.Content.Find.Execute($Find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $Replace, 2) | Out-Null
}

# Screen name corrections (Tela## references)
Replace-Text "Tela09_TelaEntrada" "Tela05_TelaEntrada"
Replace-Text "Tela14_InscricaoNaCorrida1" "Tela10_InscricaoNaCorrida1"
Replace-Text "Tela15_InscricaoNaCorrida2" "Tela11_InscricaoNaCorrida2"
Replace-Text "Tela35_pagamento" "Tela48_pagamento"
Replace-Text "Tela19_ConclusaoDePagamento" "Tela13_ConclusaoDePagamento"
Replace-Text "Tela36_pagamentoNãoAprovado" "Tela48_pagamentoNãoAprovado"

# Wording fix for the final step of the main flow
Replace-Text "Voltar para tela inicial" "Ir para tela inicial"
